# GradeBook update: quiz 3 and old work
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Letter" grade column header
$ws.Range("O1").Value = "Letter"

# Row 2 - Gary He
$ws.Range("K2").Formula = "=30/30"
$ws.Range("O2").Value = "A-"

# Row 3 - Jayden Warlum
$ws.Range("K3").Formula = "=30/30"
$ws.Range("O3").Value = "A-"

# Row 4 - Aman Kumpawat
$ws.Range("K4").Formula = "=29/30"
$ws.Range("O4").Value = "A-"

# Row 5 - Nahom Anteneh (Quiz 3 already entered)
$ws.Range("O5").Value = "A-"

# Row 6 - Oswen Martinez
$ws.Range("K6").Formula = "=0"
$ws.Range("O6").Value = "F"

# Row 7 - Kai Stephens (Quiz 3 already entered)
$ws.Range("O7").Value = "A-"

# Row 8 - James Saw
$ws.Range("K8").Formula = "=27/30"
$ws.Range("O8").Value = "A-"

# Row 9 - Edward Whitesel
$ws.Range("K9").Formula = "=29/30"
$ws.Range("O9").Value = "A-"

# Row 10 - Loren Grey
$ws.Range("K10").Formula = "=30/30"
$ws.Range("O10").Value = "B"

# Row 11 - Almas Waseem (old work: Homework 5 and Quiz 2 grades came in)
$ws.Range("I11").Formula = "=15/25"
$ws.Range("J11").Formula = "=29/35"
$ws.Range("K11").Formula = "=28/30"
$ws.Range("O11").Value = "B+"

# Row 12 - Ty Carlson
$ws.Range("K12").Formula = "=0"
$ws.Range("O12").Value = "INC"

# Row 13 - Nailyn Lopez
$ws.Range("K13").Formula = "=26/30"
$ws.Range("O13").Value = "A-"

# Row 14 - Roy Kalu (Quiz 3 already entered)
$ws.Range("O14").Value = "A-"

# Row 15 - Dylan Zeledon
$ws.Range("K15").Formula = "=32/30"
$ws.Range("O15").Value = "A-"

# View state: zoomed out, selection moved to O11
$excel.ActiveWindow.Zoom = 69
$ws.Range("O11").Select() | Out-Null
